# Add a "name" (string) field/row to the "Link" entity box, and grow the
# shape to fit the extra line of text (matches commit "Added name
# parameter to Links").

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)

# "Rectangle 19" is the "Link" table/box (Shapes collection index 10).
$sh = $s.Shapes.Item("Rectangle 19")

$tr = $sh.TextFrame.TextRange

# Paragraphs: 1 "Link", 2 "", 3 "id\tint", 4 "url\tstring", 5 "project\tFK",
# 6 "active\tboolean". Insert the new "name\tstring" row right after the
# "id\tint" row (paragraph 3), matching the XML diff ordering.
$idPara = $tr.Paragraphs(3)
[void]$idPara.InsertAfter("`rname`tstring")

# Grow the box's height (width/position stay the same) to accommodate the
# newly added line: 1973259 EMU -> 2483100 EMU (155.3747pt -> 195.5197pt).
# Shape.Height is a single-precision COM property, so nudge by a hair past
# the exact target to avoid losing the last EMU to float32 truncation.
$sh.Height = [double]195.51969146728518
